$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

foreach ($addr in @("I1", "J1")) {
    $r = $ws.Range($addr)
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4160
    $r.Font.Bold = $true
    $r.Borders.LineStyle = 1
}

$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 9
